# Commit: "add the NA's under duplicate_image_filename"
#
# Column E on the sheet is "duplicate_image_filename" (see header row E1).
# The practice trials (rows 2-5), the real trials (rows 6-13) and the
# "unique" filler rows (rows 14-21) did not yet have a value in that
# column; fill them all in with "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
